$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" - updated top-level P&L / trade-count metrics
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.8    # Current Capital
$wsSummary.Range("B4").Value = 0.9       # Total P&L $
$wsSummary.Range("B5").Value = 0.32      # Total P&L %
$wsSummary.Range("B6").Value = 57        # Total Trades
$wsSummary.Range("B7").Value = 32        # Winning Trades
$wsSummary.Range("B9").Value = 56.14     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status" - HighProbConvergence row now reflects the closed
# trade (capital/trades/P&L/win-rate all move off their zeroed defaults)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C3").Value = 100.03
$wsStatus.Range("D3").Value = 1
$wsStatus.Range("E3").Value = 0.03
$wsStatus.Range("F3").Value = 0.03
$wsStatus.Range("G3").Value = 100

# ---------------------------------------------------------------------------
# Sheet "All Trades" - trade #59 (HighProbConvergence) closes via early_exit
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")
$wsAll.Range("G60").Value = 0.6
$wsAll.Range("H60").Value = "CLOSED"
$wsAll.Range("I60").Value = 5.2632
$wsAll.Range("J60").Value = 0.03
$wsAll.Range("K60").Value = 100.03
$wsAll.Range("L60").Value = "early_exit"
$wsAll.Range("M60").Value = 0.13

# New trade #88 (momentum, OPEN) appended as row 89
$wsAll.Range("A89").Value = 88
$wsAll.Range("B89").Value = "'2026-02-18"
$wsAll.Range("C89").Value = "00:13:25"
$wsAll.Range("D89").Value = "momentum"
$wsAll.Range("E89").Value = "DOWN"
$wsAll.Range("F89").Value = 0.57
$wsAll.Range("H89").Value = "OPEN"
$wsAll.Range("I89").Value = 0
$wsAll.Range("J89").Value = 0
$wsAll.Range("K89").Value = 99.93000000000001
$wsAll.Range("M89").Value = 0
$wsAll.Range("N89").Value = 0
$wsAll.Range("O89").Value = 0
$wsAll.Range("P89").Value = 0.9
$wsAll.Range("Q89").Value = "Downward momentum: -1.980% over 10 samples"

# ---------------------------------------------------------------------------
# Sheet "momentum" - strategy-local log gains the same new trade #88 as row 19
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")
$wsMomentum.Range("A19").Value = 88
$wsMomentum.Range("B19").Value = "'2026-02-18"
$wsMomentum.Range("C19").Value = "00:13:25"
$wsMomentum.Range("D19").Value = "momentum"
$wsMomentum.Range("E19").Value = "DOWN"
$wsMomentum.Range("F19").Value = 0.57
$wsMomentum.Range("H19").Value = "OPEN"
$wsMomentum.Range("I19").Value = 0
$wsMomentum.Range("J19").Value = 0
$wsMomentum.Range("K19").Value = 99.93000000000001
$wsMomentum.Range("L19").Value = 0
$wsMomentum.Range("M19").Value = 0
$wsMomentum.Range("N19").Value = 0.9
$wsMomentum.Range("O19").Value = "Downward momentum: -1.980% over 10 samples"
$wsMomentum.Range("Q19").Value = 0

# ---------------------------------------------------------------------------
# Sheet "HighProbConvergence" - strategy-local log: trade #59 closes
# ---------------------------------------------------------------------------
$wsHPC = $wb.Worksheets.Item("HighProbConvergence")
$wsHPC.Range("G2").Value = 0.6
$wsHPC.Range("H2").Value = "CLOSED"
$wsHPC.Range("I2").Value = 5.2632
$wsHPC.Range("J2").Value = 0.03
$wsHPC.Range("K2").Value = 100.03
$wsHPC.Range("P2").Value = "early_exit"
$wsHPC.Range("Q2").Value = 0.13
